$wb = $excel.ActiveWorkbook

# --- Data edit: "lethality" sheet gets an extra lethal-sample scenario row ---
# Existing row 4 (scenario 3) drops from 1000 lethal samples/yr to 500,
# and a new row 5 (scenario 4) is added carrying the 1000/yr values that
# row 4 used to hold.
$wsLethality = $wb.Worksheets.Item("lethality")
$wsLethality.Range("B4:F4").Value = 500
$wsLethality.Range("A5").Value = 4
$wsLethality.Range("B5:F5").Value = 1000

# --- View/selection state to match the author's final on-screen state ---
# Demography sheet: selection moved to A5:K7 (anchor A5)
$wsDemography = $wb.Worksheets.Item("demography")
$wsDemography.Activate() | Out-Null
$wsDemography.Range("A5:K7").Select() | Out-Null

# Finally, the lethality sheet is the active tab when the workbook was saved.
$wsLethality.Activate() | Out-Null
